$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A11").NumberFormat = "@"

# Update existing rows 2-8 (amount_value / count pairs)
$ws.Range("A2").Value = "220.00"
$ws.Range("B2").Value = 109

$ws.Range("A3").Value = "440.00"
$ws.Range("B3").Value = 24

$ws.Range("A4").Value = "320.00"
$ws.Range("B4").Value = 1

$ws.Range("A5").Value = "178.71"
$ws.Range("B5").Value = 1

$ws.Range("A6").Value = "1320.00"
$ws.Range("B6").Value = 3

$ws.Range("A7").Value = "660.00"
$ws.Range("B7").Value = 5

$ws.Range("A8").Value = "120.00"
$ws.Range("B8").Value = 1

# New rows 9-11
$ws.Range("A9").Value = "300.00"
$ws.Range("B9").Value = 1

$ws.Range("A10").Value = "2640.00"
$ws.Range("B10").Value = 1

$ws.Range("A11").Value = "1100.00"
$ws.Range("B11").Value = 1
